$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row 84 by copying the formatting/structure of row 83 (A:R),
# which will be overwritten below with the correct shifted values.
$ws.Range("A83:R83").Copy($ws.Range("A84:R84"))

# Update Fecha (D), Volumen (J), Precio minimo (K), Precio maximo (L),
# Precio promedio ponderado (M) and Precio $/Kg (P) for rows 62-84
# to reflect the corrected weekly data (values shift down one row,
# with a new record appended at the end).
$ws.Range("D62").Value = 44468
$ws.Range("J62").Value = 600
$ws.Range("K62").Value = 8000
$ws.Range("L62").Value = 9000
$ws.Range("M62").Value = 8500
$ws.Range("P62").Value = 142

$ws.Range("D63").Value = 44312
$ws.Range("J63").Value = 600
$ws.Range("K63").Value = 8000
$ws.Range("L63").Value = 9000
$ws.Range("M63").Value = 8500
$ws.Range("P63").Value = 142

$ws.Range("D64").Value = 44386
$ws.Range("J64").Value = 560
$ws.Range("K64").Value = 11000
$ws.Range("L64").Value = 12000
$ws.Range("M64").Value = 11500
$ws.Range("P64").Value = 192

$ws.Range("D65").Value = 44463
$ws.Range("J65").Value = 600
$ws.Range("K65").Value = 9000
$ws.Range("L65").Value = 10000
$ws.Range("M65").Value = 9500
$ws.Range("P65").Value = 158

$ws.Range("D66").Value = 44251
$ws.Range("J66").Value = 600
$ws.Range("K66").Value = 8000
$ws.Range("L66").Value = 9000
$ws.Range("M66").Value = 8500
$ws.Range("P66").Value = 142

$ws.Range("D67").Value = 44371
$ws.Range("J67").Value = 560
$ws.Range("K67").Value = 13000
$ws.Range("L67").Value = 14000
$ws.Range("M67").Value = 13500
$ws.Range("P67").Value = 225

$ws.Range("D68").Value = 44316
$ws.Range("J68").Value = 520
$ws.Range("K68").Value = 8000
$ws.Range("L68").Value = 9000
$ws.Range("M68").Value = 8500
$ws.Range("P68").Value = 142

$ws.Range("D69").Value = 44279
$ws.Range("J69").Value = 600
$ws.Range("K69").Value = 8000
$ws.Range("L69").Value = 9000
$ws.Range("M69").Value = 8500
$ws.Range("P69").Value = 142

$ws.Range("D70").Value = 44397
$ws.Range("J70").Value = 560
$ws.Range("K70").Value = 12000
$ws.Range("L70").Value = 12500
$ws.Range("M70").Value = 12250
$ws.Range("P70").Value = 204

$ws.Range("D71").Value = 44372
$ws.Range("J71").Value = 560
$ws.Range("K71").Value = 13000
$ws.Range("L71").Value = 14000
$ws.Range("M71").Value = 13500
$ws.Range("P71").Value = 225

$ws.Range("D72").Value = 44286
$ws.Range("J72").Value = 600
$ws.Range("K72").Value = 8000
$ws.Range("L72").Value = 9000
$ws.Range("M72").Value = 8500
$ws.Range("P72").Value = 142

$ws.Range("D73").Value = 44356
$ws.Range("J73").Value = 600
$ws.Range("K73").Value = 12000
$ws.Range("L73").Value = 13000
$ws.Range("M73").Value = 12500
$ws.Range("P73").Value = 208

$ws.Range("D74").Value = 44160
$ws.Range("J74").Value = 700
$ws.Range("K74").Value = 9000
$ws.Range("L74").Value = 10000
$ws.Range("M74").Value = 9500
$ws.Range("P74").Value = 158

$ws.Range("D75").Value = 44351
$ws.Range("J75").Value = 520
$ws.Range("K75").Value = 11500
$ws.Range("L75").Value = 12000
$ws.Range("M75").Value = 11750
$ws.Range("P75").Value = 196

$ws.Range("D76").Value = 44365
$ws.Range("J76").Value = 520
$ws.Range("K76").Value = 13000
$ws.Range("L76").Value = 14000
$ws.Range("M76").Value = 13500
$ws.Range("P76").Value = 225

$ws.Range("D77").Value = 44306
$ws.Range("J77").Value = 500
$ws.Range("K77").Value = 8500
$ws.Range("L77").Value = 9000
$ws.Range("M77").Value = 8750
$ws.Range("P77").Value = 146

$ws.Range("D78").Value = 44357
$ws.Range("J78").Value = 520
$ws.Range("K78").Value = 12000
$ws.Range("L78").Value = 12500
$ws.Range("M78").Value = 12250
$ws.Range("P78").Value = 204

$ws.Range("D79").Value = 44162
$ws.Range("J79").Value = 600
$ws.Range("K79").Value = 9000
$ws.Range("L79").Value = 10000
$ws.Range("M79").Value = 9500
$ws.Range("P79").Value = 158

$ws.Range("D80").Value = 44244
$ws.Range("J80").Value = 600
$ws.Range("K80").Value = 8000
$ws.Range("L80").Value = 9000
$ws.Range("M80").Value = 8500
$ws.Range("P80").Value = 142

$ws.Range("D81").Value = 44239
$ws.Range("J81").Value = 800
$ws.Range("K81").Value = 8000
$ws.Range("L81").Value = 9000
$ws.Range("M81").Value = 8500
$ws.Range("P81").Value = 142

$ws.Range("D82").Value = 44376
$ws.Range("J82").Value = 520
$ws.Range("K82").Value = 12000
$ws.Range("L82").Value = 13000
$ws.Range("M82").Value = 12500
$ws.Range("P82").Value = 208

$ws.Range("D83").Value = 44358
$ws.Range("J83").Value = 540
$ws.Range("K83").Value = 11500
$ws.Range("L83").Value = 12000
$ws.Range("M83").Value = 11750
$ws.Range("P83").Value = 196

$ws.Range("D84").Value = 44323
$ws.Range("J84").Value = 500
$ws.Range("K84").Value = 8000
$ws.Range("L84").Value = 9000
$ws.Range("M84").Value = 8500
$ws.Range("P84").Value = 142
